$d = $word.ActiveDocument

# The document currently ends with a paragraph that contains only the
# "_GoBack" bookmark (no visible text). Word's Range.InsertXML merges the
# *last* <w:p> of an inserted fragment into the paragraph at the insertion
# point (keeping whatever already followed the insertion point - here, the
# bookmark - inside that same merged paragraph), while every earlier <w:p>
# in the fragment becomes its own brand-new paragraph. We exploit that by
# ending the first fragment with the "Stage 2 - Jungle" paragraph, so it
# ends up sharing a paragraph with the bookmark, exactly like the target
# diff. The remaining paragraphs (after "Stage 2 - Jungle") are inserted
# with a second call positioned right after that merged paragraph.

$wordmlNamespace = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-WordXmlPackage([string]$bodyFragment) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
      '<w:document ' + $wordmlNamespace + '>' +
      '<w:body>' + $bodyFragment + '</w:body>' +
      '</w:document>' +
      '</pkg:xmlData>' +
      '</pkg:part>' +
      '</pkg:package>'
}

# --- Part 1: everything from "Stage 1 - In the clouds" through the
#     "Stage 2 - Jungle" heading (inclusive). The final paragraph in this
#     fragment merges into the bookmark's paragraph.
$bodyFragment1 = @'
<w:p/><w:p/><w:p/><w:p/><w:p><w:r><w:t>Stage 1 – In the clouds</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>{image}</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">The level system for this game includes two levels which contain three stages. An idea for this system would be to have different locations in the jungle. One of the stages could be above the jungle’s trees or in the clouds. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Low fidelity designs have been created to visualize these ideas. Including these different locations will allow for different scenarios. Thus, giving the player the opportunity to use the energy drinks and the power ups to move onto the next stage. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">In this example stage, the player would be able to consume two of the six types of energy drinks; Grape (jump higher) and Apple (Invisibility). These drinks correspond with this level the most because the whole idea would be to jump onto clouds and not fall </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>and also</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> avoiding enemies in the air. </w:t></w:r></w:p><w:p/><w:p><w:r><w:br w:type="page"/></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>Stage 2 - Jungle</w:t></w:r></w:p>
'@

$bm = $d.Bookmarks.Item("_GoBack")
$insertPos = $bm.Range.Start
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertXML((New-WordXmlPackage $bodyFragment1)) | Out-Null

# --- Part 2: everything after the "Stage 2 - Jungle" paragraph, through to
#     the end of the document body (before the sectPr).
$bodyFragment2 = @'
<w:p/><w:p><w:r><w:t>{image}</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr></w:p><w:p><w:r><w:t xml:space="preserve">Another location amongst the stages would be the jungle floor. In this stage, the player would be able to consume two of the six types of energy drinks; Blueberry (throw rocks) and Banana (Super strength). These energy drinks would be the most appropriate for this this stage, as this level will contain a lot of trying to defeat enemies </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>and also</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> possibly moving heavy objects. The only way to get through these obstacles will be to consume the drinks.</w:t></w:r></w:p><w:p/><w:p/><w:p/><w:p/><w:p><w:r><w:t>Stage 3 – Caves / Underground</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>{image}</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">The final location choice for the game would be underground or in some sort of cave system. The main idea of the different locations is trying to create various scenarios all within the jungle but also having </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>a very obvious</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> change in scenery. This will increase the players interest, as it’s not all in the same location.</w:t></w:r></w:p><w:p><w:r><w:t>In this stage, the player would be able to consume two of the six types of energy drinks; Carrots (Seeing in the dark) and Kiwi (Minimize). These ingredients are very appropriate to this stage because there is the possibility that the caves might be too dark, and the player will have to find the Carrot to be able to see. Or there may be a small hole in which the player can only get through, if they find the Kiwi.</w:t></w:r></w:p><w:p/>
'@

$bm2 = $d.Bookmarks.Item("_GoBack")
$stage2Paragraph = $bm2.Range.Paragraphs.First
$afterPos = $stage2Paragraph.Range.End
$afterRange = $d.Range($afterPos, $afterPos)
$afterRange.InsertXML((New-WordXmlPackage $bodyFragment2)) | Out-Null
